# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections described by the commit diff to the
# "Leve Profit" calculation columns (H..N) across all 8 class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12136210
$ws.Range("I6").Value = 3560016.8
$ws.Range("J6").Value = 25000500
$ws.Range("K6").Value = 10680050.4
$ws.Range("L6").Value = 75001500
$ws.Range("M6").Value = -10679938.4
$ws.Range("N6").Value = -75001724

$ws.Range("H9").Value = 183.33333
$ws.Range("I9").Value = 100
$ws.Range("J9").Value = 350
$ws.Range("K9").Value = 100
$ws.Range("L9").Value = 350
$ws.Range("M9").Value = 69
$ws.Range("N9").Value = -688

$ws.Range("H12").Value = 175.5
$ws.Range("I12").Value = 175.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 175.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -5.5
$ws.Range("N12").ClearContents()

$ws.Range("H21").Value = 3136.5715
$ws.Range("I21").Value = 1986.6666
$ws.Range("J21").Value = 3999
$ws.Range("K21").Value = 1986.6666
$ws.Range("L21").Value = 3999
$ws.Range("M21").Value = -1518.6666
$ws.Range("N21").Value = -4935

$ws.Range("H23").Value = 3136.5715
$ws.Range("I23").Value = 1986.6666
$ws.Range("J23").Value = 3999
$ws.Range("K23").Value = 1986.6666
$ws.Range("L23").Value = 3999
$ws.Range("M23").Value = -1752.6666
$ws.Range("N23").Value = -4467

$ws.Range("H29").Value = 663
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 717.3333
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 2151.9999
$ws.Range("M29").Value = -1219
$ws.Range("N29").Value = -2713.9999

$ws.Range("H38").Value = 1862.7391
$ws.Range("I38").Value = 208.7
$ws.Range("J38").Value = 3135.077
$ws.Range("K38").Value = 626.0999999999999
$ws.Range("L38").Value = 9405.231
$ws.Range("M38").Value = -254.0999999999999
$ws.Range("N38").Value = -10149.231

$ws.Range("H40").Value = 1838.5834
$ws.Range("I40").Value = 1912.2
$ws.Range("J40").Value = 1786
$ws.Range("K40").Value = 1912.2
$ws.Range("L40").Value = 1786
$ws.Range("M40").Value = -1737.2
$ws.Range("N40").Value = -2136

$ws.Range("H42").Value = 259.6
$ws.Range("I42").Value = 206.875
$ws.Range("J42").Value = 319.85715
$ws.Range("K42").Value = 620.625
$ws.Range("L42").Value = 959.5714499999999
$ws.Range("M42").Value = -390.625
$ws.Range("N42").Value = -1419.57145

$ws.Range("H43").Value = 569.8484999999999
$ws.Range("I43").Value = 490.6
$ws.Range("J43").Value = 584
$ws.Range("K43").Value = 490.6
$ws.Range("L43").Value = 584
$ws.Range("M43").Value = -421.6
$ws.Range("N43").Value = -722

$ws.Range("H58").Value = 561.96
$ws.Range("I58").Value = 123.818184
$ws.Range("J58").Value = 906.2143
$ws.Range("K58").Value = 371.454552
$ws.Range("L58").Value = 2718.6429
$ws.Range("M58").Value = -221.454552
$ws.Range("N58").Value = -3018.6429

$ws.Range("H62").Value = 1603.4
$ws.Range("I62").Value = 1800
$ws.Range("J62").Value = 1406.8
$ws.Range("K62").Value = 1800
$ws.Range("L62").Value = 1406.8
$ws.Range("M62").Value = -1176
$ws.Range("N62").Value = -2654.8

$ws.Range("H65").Value = 1603.4
$ws.Range("I65").Value = 1800
$ws.Range("J65").Value = 1406.8
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 7034
$ws.Range("M65").Value = -5880
$ws.Range("N65").Value = -13274

$ws.Range("H87").Value = 11864.383
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 11864.383
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 11864.383
$ws.Range("N87").Value = -14360.383

$ws.Range("H90").Value = 11864.383
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 11864.383
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 35593.149
$ws.Range("N90").Value = -48073.149

$ws.Range("H107").Value = 658
$ws.Range("I107").Value = 802.8095
$ws.Range("J107").Value = 223.57143
$ws.Range("K107").Value = 802.8095
$ws.Range("L107").Value = 223.57143
$ws.Range("M107").Value = 1117.1905
$ws.Range("N107").Value = -4063.57143

$ws.Range("H113").Value = 2025.75
$ws.Range("I113").Value = 2206.111
$ws.Range("J113").Value = 1701.1
$ws.Range("K113").Value = 2206.111
$ws.Range("L113").Value = 1701.1
$ws.Range("M113").Value = 1047.889
$ws.Range("N113").Value = -8209.1

$ws.Range("H132").Value = 33492.594
$ws.Range("I132").Value = 38957.89
$ws.Range("J132").Value = 3980
$ws.Range("K132").Value = 116873.67
$ws.Range("L132").Value = 11940
$ws.Range("M132").Value = -114343.67
$ws.Range("N132").Value = -17000

$ws.Range("H138").Value = 2929.6428
$ws.Range("I138").Value = 1270.7778
$ws.Range("J138").Value = 3971.2559
$ws.Range("K138").Value = 3812.3334
$ws.Range("L138").Value = 11913.7677
$ws.Range("M138").Value = 1327.6666
$ws.Range("N138").Value = -22193.7677

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -388

$ws.Range("H63").Value = 3832.5
$ws.Range("I63").Value = 2165.5
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 2165.5
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1479.5
$ws.Range("N63").Value = -9372

$ws.Range("H66").Value = 3832.5
$ws.Range("I66").Value = 2165.5
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 10827.5
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -7395.5
$ws.Range("N66").Value = -46864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -385

$ws.Range("H15").Value = 2000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2000
$ws.Range("N15").Value = -2454

$ws.Range("H20").Value = 3554.2188
$ws.Range("I20").Value = 2945.818
$ws.Range("J20").Value = 4892.7
$ws.Range("K20").Value = 2945.818
$ws.Range("L20").Value = 4892.7
$ws.Range("M20").Value = -2698.818
$ws.Range("N20").Value = -5386.7

$ws.Range("H82").Value = 1999.5
$ws.Range("I82").Value = 1999.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1999.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1616.5
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 1999.5
$ws.Range("I85").Value = 1999.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1999.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -673.5
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6135.148
$ws.Range("I31").Value = 6336.913
$ws.Range("J31").Value = 4975
$ws.Range("K31").Value = 6336.913
$ws.Range("L31").Value = 4975
$ws.Range("M31").Value = -6041.913
$ws.Range("N31").Value = -5565

$ws.Range("H34").Value = 6135.148
$ws.Range("I34").Value = 6336.913
$ws.Range("J34").Value = 4975
$ws.Range("K34").Value = 6336.913
$ws.Range("L34").Value = 4975
$ws.Range("M34").Value = -6134.913
$ws.Range("N34").Value = -5379

$ws.Range("H99").Value = 284751.66
$ws.Range("I99").Value = 328340.56
$ws.Range("J99").Value = 1423.75
$ws.Range("K99").Value = 328340.56
$ws.Range("L99").Value = 1423.75
$ws.Range("M99").Value = -326842.56
$ws.Range("N99").Value = -4419.75

$ws.Range("H100").Value = 32640
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 32640
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 32640
$ws.Range("N100").Value = -34804

$ws.Range("H126").Value = 284751.66
$ws.Range("I126").Value = 328340.56
$ws.Range("J126").Value = 1423.75
$ws.Range("K126").Value = 985021.6799999999
$ws.Range("L126").Value = 4271.25
$ws.Range("M126").Value = -982551.6799999999
$ws.Range("N126").Value = -9211.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 857.35
$ws.Range("I5").Value = 895.3889
$ws.Range("J5").Value = 515
$ws.Range("K5").Value = 2686.1667
$ws.Range("L5").Value = 1545
$ws.Range("M5").Value = -2574.1667
$ws.Range("N5").Value = -1769

$ws.Range("H39").Value = 2077.2727
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2077.2727
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 6231.8181
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -6819.8181

$ws.Range("H51").Value = 1814.2858
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1814.2858
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5442.857400000001
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -6362.857400000001

$ws.Range("H86").Value = 1995
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1995
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5985
$ws.Range("N86").Value = -8357

$ws.Range("H89").Value = 1995
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1995
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 17955
$ws.Range("N89").Value = -29811

$ws.Range("H113").Value = 797.45
$ws.Range("I113").Value = 650.6667
$ws.Range("J113").Value = 829.6707
$ws.Range("K113").Value = 1952.0001
$ws.Range("L113").Value = 2489.0121
$ws.Range("M113").Value = 217.9999
$ws.Range("N113").Value = -6829.0121

$ws.Range("H131").Value = 6411233.5
$ws.Range("I131").Value = 1185.7142
$ws.Range("J131").Value = 7043210
$ws.Range("K131").Value = 3557.1426
$ws.Range("L131").Value = 21129630
$ws.Range("M131").Value = 1482.8574
$ws.Range("N131").Value = -21139710

$ws.Range("H135").Value = 857.35
$ws.Range("I135").Value = 895.3889
$ws.Range("J135").Value = 515
$ws.Range("K135").Value = 8058.5001
$ws.Range("L135").Value = 4635
$ws.Range("M135").Value = -5523.5001
$ws.Range("N135").Value = -9705

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 23824626
$ws.Range("I126").Value = 20842.4
$ws.Range("J126").Value = 83334090
$ws.Range("K126").Value = 62527.2
$ws.Range("L126").Value = 250002270
$ws.Range("M126").Value = -60057.2
$ws.Range("N126").Value = -250007210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1531.05
$ws.Range("I40").Value = 1431.9584
$ws.Range("J40").Value = 1679.6875
$ws.Range("K40").Value = 1431.9584
$ws.Range("L40").Value = 1679.6875
$ws.Range("M40").Value = -1295.9584
$ws.Range("N40").Value = -1951.6875

$ws.Range("H122").Value = 2423.6191
$ws.Range("I122").Value = 2241.4167
$ws.Range("J122").Value = 2666.5557
$ws.Range("K122").Value = 6724.250100000001
$ws.Range("L122").Value = 7999.6671
$ws.Range("M122").Value = -4274.250100000001
$ws.Range("N122").Value = -12899.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 33609.438
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 33609.438
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 33609.438
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -43409.438

